# Update the "Raju Ahamed" sheet with the 18.06.19 Today Sales data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# Update the date headers (both tables on this sheet) from 17.06.19 to 18.06.19
$ws.Range("A4").Value = "Date: 18.06.19"
$ws.Range("A31").Value = "Date: 18.06.19"

# Update the Qty. (E) column figures for the first table (rows 6-11)
$ws.Range("E6").Value = 44
$ws.Range("E7").Value = 141
$ws.Range("E8").Value = 213
$ws.Range("E9").Value = 24
$ws.Range("E10").Value = $null

# Update the Qty. (E) column figures for the second (duplicate) table (rows 33-38)
$ws.Range("E33").Value = 44
$ws.Range("E34").Value = 141
$ws.Range("E35").Value = 213
$ws.Range("E36").Value = 24
$ws.Range("E37").Value = $null
